# Risoluzione di piccoli errori e debugging
# Updates the "Schedulazione" sheet rows 4-12 (scheduling heuristic output)
# with the re-run results: commessa order, setup/processing minutes,
# timestamps, compatible-machine lists, knife counts, vehicle refs and
# lateness (R column) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 251752
$ws.Range("C4").Value = 17
$ws.Range("D4").Value = 0
$ws.Range("F4").Value = "2025-06-04 12:17:00"
$ws.Range("G4").Value = "2025-06-04 12:17:00"
$ws.Range("H4").Value = "2025-06-04 12:17:00"
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R6 ;R9"
$ws.Range("L4").Value = 3
$ws.Range("N4").Value = 39846
$ws.Range("P4").Value = 39846
$ws.Range("Q4").Value = "2025-05-20 00:00:00"
$ws.Range("R4").Value = -0.5118055555555555
$ws.Range("A5").Value = 251218
$ws.Range("D5").Value = 96.90140845070422
$ws.Range("E5").Value = "2025-06-04 12:17:00"
$ws.Range("F5").Value = "2025-06-04 12:38:00"
$ws.Range("G5").Value = "2025-06-04 12:38:00"
$ws.Range("H5").Value = "2025-06-04 14:14:54"
$ws.Range("I5").Value = 6880
$ws.Range("K5").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R9"
$ws.Range("L5").Value = 6
$ws.Range("N5").Value = 39885
$ws.Range("P5").Value = 39885
$ws.Range("Q5").Value = "2025-05-09 00:00:00"
$ws.Range("R5").Value = 0
$ws.Range("A6").Value = 251565
$ws.Range("C6").Value = 38
$ws.Range("D6").Value = 176.7464788732394
$ws.Range("E6").Value = "2025-06-04 14:14:54"
$ws.Range("H6").Value = "2025-06-05 09:49:38"
$ws.Range("I6").Value = 12549
$ws.Range("K6").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$ws.Range("L6").Value = 2
$ws.Range("Q6").Value = "2025-06-10 00:00:00"
$ws.Range("R6").Value = -0.409477699525463
$ws.Range("S6").Value = 1
$ws.Range("A7").Value = 251500
$ws.Range("D7").Value = 139.3802816901408
$ws.Range("E7").Value = "2025-06-05 09:49:38"
$ws.Range("F7").Value = "2025-06-05 10:08:38"
$ws.Range("G7").Value = "2025-06-05 10:08:38"
$ws.Range("I7").Value = 9896
$ws.Range("K7").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R6 ;R9"
$ws.Range("L7").Value = 4
$ws.Range("Q7").Value = "2025-05-26 00:00:00"
$ws.Range("S7").Value = 2
$ws.Range("A8").Value = 251070
$ws.Range("C8").Value = 36.5
$ws.Range("F8").Value = "2025-06-05 07:36:30"
$ws.Range("G8").Value = "2025-06-05 07:36:30"
$ws.Range("H8").Value = "2025-06-05 07:36:30"
$ws.Range("K8").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R6 ;R9"
$ws.Range("L8").Value = 6
$ws.Range("N8").Value = 39885
$ws.Range("P8").Value = 39885
$ws.Range("Q8").Value = "2025-03-28 00:00:00"
$ws.Range("R8").Value = -0.3170138888888889
$ws.Range("S8").Value = 2
$ws.Range("A9").Value = 251773
$ws.Range("C9").Value = 32.5
$ws.Range("E9").Value = "2025-06-05 07:36:30"
$ws.Range("F9").Value = "2025-06-05 08:09:00"
$ws.Range("G9").Value = "2025-06-05 08:09:00"
$ws.Range("H9").Value = "2025-06-05 08:09:00"
$ws.Range("K9").Value = "CASON ;R6"
$ws.Range("L9").Value = 7
$ws.Range("N9").Value = 39874
$ws.Range("P9").Value = 39874
$ws.Range("Q9").Value = "2025-05-25 00:00:00"
$ws.Range("R9").Value = -0.3395833333333333
$ws.Range("S9").Value = 1
$ws.Range("A10").Value = 251180
$ws.Range("C10").Value = 36.5
$ws.Range("E10").Value = "2025-06-05 08:09:00"
$ws.Range("F10").Value = "2025-06-05 08:45:30"
$ws.Range("G10").Value = "2025-06-05 08:45:30"
$ws.Range("H10").Value = "2025-06-05 08:45:30"
$ws.Range("K10").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$ws.Range("L10").Value = 4
$ws.Range("N10").Value = "39887 (esterno)"
$ws.Range("P10").Value = 39887
$ws.Range("Q10").Value = "2025-05-20 00:00:00"
$ws.Range("R10").Value = -16.36493055555556
$ws.Range("S10").Value = 7
$ws.Range("C11").Value = 42.5
$ws.Range("E11").Value = "2025-06-05 08:45:30"
$ws.Range("F11").Value = "2025-06-05 09:28:00"
$ws.Range("G11").Value = "2025-06-05 09:28:00"
$ws.Range("H11").Value = "2025-06-05 14:49:42"
$ws.Range("R11").Value = -10.6178535353588
$ws.Range("E12").Value = "2025-06-05 14:49:42"
$ws.Range("F12").Value = "2025-06-06 07:34:12"
$ws.Range("G12").Value = "2025-06-06 07:34:12"
$ws.Range("H12").Value = "2025-06-06 14:56:00"
